# story_aglina_1_1.xlsx carries a 4-column glossary table:
#   A1:D1 = zh_CN / ja_JP / en_US / ko_KR   (language headers)
#   A2:D2 = the matching translated blurb for each language
#
# Row 2 previously only had the zh_CN and ko_KR blurbs filled in (with the
# ja_JP / en_US cells temporarily holding a duplicate of the Chinese text).
# This change fills in the real Japanese and English translations.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# B2 = ja_JP translation
$ws.Range("B2").Value = "受取人が書かれていない荷物を見つけたアンジェリーナはクロワッサンと共に、その持ち主を探し始める。`n"

# C2 = en_US translation
$ws.Range("C2").Value = "Angelina finds an unaddressed package. With Croissant's help, she searches for its owner.`n"

# D2 = ko_KR translation (unchanged text, re-applied so the cell is
# consistently refreshed along with the rest of the row)
$ws.Range("D2").Value = "안젤리나는 수신인이 적혀 있지 않은 소포를 발견했다. 그녀는 크루아상의 도움을 받아 소포의 주인을 찾아 나서게 되는데……`n"
